$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Category" -> "Description"
$ws.Range("A1").Value = "Description"

# Row 2: Groceries/250/2025-07-20 -> Rent/300/2025-07-24
$ws.Range("A2").Value = "Rent"
$ws.Range("B2").Value = 300
$ws.Range("C2").Value = 45862

# New row 3: Cat Food/70/2025-07-23
$ws.Range("A3").Value = "Cat Food"
$ws.Range("B3").Value = 70
$ws.Range("C3").Value = 45861

# Copy C2's date format onto the new C3 cell so it keeps the same
# short-date number format style as the rest of the Date column.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
